$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.171.63'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.302.77'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.88'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.53'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.26%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.501'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.77'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.73%  '
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.27'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.90%  '
$ws.Range("E13").Value = '  +2.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.03'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.31%  '
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").Value = '2.659.67'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = '2.307.31'
$ws.Range("E17").Value = '  -0.70%  '
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").Value = '43.054.92'
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.61'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.07'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '237.20'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.02'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.45'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.58'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.88'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.10'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.02%  '
$ws.Range("E31").Value = '  -5.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.15'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.69'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.96'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '16.98'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.92%  '
$ws.Range("E37").Value = '  -1.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0695'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.81%  '
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.77'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.46%  '
$ws.Range("E42").Value = '  -1.03%  '
$ws.Range("E43").Value = '  -2.23%  '
$ws.Range("D44").Value = '1.982.80'
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("E45").Value = '  -1.31%  '
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.67'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.87'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.33%  '
$ws.Range("D49").Value = '2.527.38'
$ws.Range("E49").Value = '  -0.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.31'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.70%  '
$ws.Range("E51").Value = '  -4.16%  '
